$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-10"

# Update the header label for the first data column
$ws.Range("B1").Value = "March 2022 (through March 10)"

# Row 3 - Austin
$ws.Range("Q3").Value = 2

# Row 7 - South Shore
$ws.Range("E7").Value = 3

# Row 12 - Englewood
$ws.Range("B12").Value = 5
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 2
$ws.Range("N12").Value = 2

# Row 15 - Humboldt Park
$ws.Range("Q15").Value = 2

# Row 17 - Auburn Gresham
$ws.Range("T17").Value = 1

# Row 61 - Fuller Park
$ws.Range("E61").Value = 1
$ws.Range("Q61").Value = 1

# Row 84 - South Deering
$ws.Range("B84").Value = 1
